# Corrections in raw data
# Updates the raw conductivity measurements (column C) on the
# "ALD TiN_conductivity" sheet, clears the now-unneeded explicit
# number-format / highlight styling that accompanied the old values,
# and leaves the workbook with sheet 2 ("ALD TiN_conductivity") as the
# active / selected sheet with C36 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Corrected raw data values for column C (rows 2-36)
# ---------------------------------------------------------------------
$corrections = @(
    @{ Row = 2;  Value = 382.59999299999998 },
    @{ Row = 3;  Value = 380.86066199999999 },
    @{ Row = 4;  Value = 360.766051 },
    @{ Row = 5;  Value = 370.18898200000001 },
    @{ Row = 6;  Value = 375.33335499999998 },
    @{ Row = 7;  Value = 393.972579 },
    @{ Row = 8;  Value = 322.33842600000003 },
    @{ Row = 9;  Value = 371.64198699999997 },
    @{ Row = 10; Value = 349.65438 },
    @{ Row = 11; Value = 376.87962499999998 },
    @{ Row = 12; Value = 318.32039400000002 },
    @{ Row = 13; Value = 347.73478 },
    @{ Row = 14; Value = 353.41118399999999 },
    @{ Row = 15; Value = 301.34893299999999 },
    @{ Row = 16; Value = 298.35282699999999 },
    @{ Row = 17; Value = 327.83011499999998 },
    @{ Row = 18; Value = 323.36190099999999 },
    @{ Row = 19; Value = 281.89321899999999 },
    @{ Row = 20; Value = 275.018033 },
    @{ Row = 21; Value = 283.72517699999997 },
    @{ Row = 22; Value = 361.77751899999998 },
    @{ Row = 23; Value = 328.60957000000002 },
    @{ Row = 24; Value = 332.90978000000001 },
    @{ Row = 25; Value = 279.81870800000002 },
    @{ Row = 26; Value = 335.69217200000003 },
    @{ Row = 27; Value = 378.53232400000002 },
    @{ Row = 28; Value = 348.51854100000003 },
    @{ Row = 29; Value = 342.366196 },
    @{ Row = 30; Value = 329.00780500000002 },
    @{ Row = 31; Value = 362.90587499999998 },
    @{ Row = 32; Value = 311.83471500000002 },
    @{ Row = 33; Value = 283.00771300000002 },
    @{ Row = 34; Value = 345.32294899999999 },
    @{ Row = 35; Value = 353.65588200000002 },
    @{ Row = 36; Value = 347.76735200000002 }
)

foreach ($item in $corrections) {
    $ws2.Range("C" + $item.Row).Value = $item.Value
}

# ---------------------------------------------------------------------
# 2. Drop the explicit "#,##0.000" number format that used to mark these
#    raw-value cells - the corrected sheet no longer special-cases them.
# ---------------------------------------------------------------------
$clearFormatRows = @(2,4,8,9,10,11,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36)
foreach ($r in $clearFormatRows) {
    $ws2.Range("C" + $r).ClearFormats()
}

# Row 21 used to carry a whole-row custom format (fill + bold-less
# highlight on every cell plus a distinct number format on C21) - strip
# all of it in one go so every cell in the row reverts to the default
# style.
$ws2.Rows(21).ClearFormats()

# The highlight fill that used to flag rows 28-36 in column A is gone too.
$ws2.Range("A28:A36").Interior.Pattern = -4142

# ---------------------------------------------------------------------
# 3. View state: sheet 2 becomes the active/selected sheet (was sheet 1)
#    with C36 selected, scrolled near the bottom of the table.
# ---------------------------------------------------------------------
$ws1.Range("F4").Select()
$ws1.Activate()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("C36").Select()
